$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value = 172
$ws.Range("E10").Value = 735
$ws.Range("F10").Value = 408
$ws.Range("H10").Value = 503
$ws.Range("E11").Value = 490
$ws.Range("F11").Value = 274
$ws.Range("H11").Value = 339
$ws.Range("E12").Value = 748
$ws.Range("F12").Value = 442
$ws.Range("H12").Value = 528
$ws.Range("E13").Value = 175
$ws.Range("F13").Value = 98
$ws.Range("H13").Value = 132
$ws.Range("E14").Value = 148
$ws.Range("E15").Value = 208
$ws.Range("E16").Value = 239
$ws.Range("F16").Value = 136
$ws.Range("H16").Value = 184
$ws.Range("E17").Value = 127
$ws.Range("E20").Value = 101
$ws.Range("F20").Value = 47
$ws.Range("H20").Value = 84
$ws.Range("E21").Value = 153
$ws.Range("E22").Value = 203
$ws.Range("F22").Value = 112
$ws.Range("H22").Value = 154
$ws.Range("E23").Value = 233
$ws.Range("F23").Value = 119
$ws.Range("H23").Value = 171
$ws.Range("E24").Value = 283
$ws.Range("F24").Value = 164
$ws.Range("H24").Value = 194
$ws.Range("E25").Value = 344
$ws.Range("E26").Value = 217
$ws.Range("F26").Value = 126
$ws.Range("H26").Value = 151
$ws.Range("E27").Value = 400
$ws.Range("F27").Value = 218
$ws.Range("H27").Value = 300
$ws.Range("E28").Value = 235
$ws.Range("F28").Value = 115
$ws.Range("H28").Value = 167
$ws.Range("E29").Value = 199
$ws.Range("E30").Value = 269
$ws.Range("F30").Value = 165
$ws.Range("H30").Value = 218
$ws.Range("E31").Value = 85
$ws.Range("E32").Value = 224
$ws.Range("F32").Value = 142
$ws.Range("H32").Value = 180
$ws.Range("E33").Value = 347
$ws.Range("F33").Value = 183
$ws.Range("H33").Value = 274
$ws.Range("E34").Value = 265
$ws.Range("F34").Value = 184
$ws.Range("H34").Value = 222
$ws.Range("E35").Value = 191
$ws.Range("F35").Value = 130
$ws.Range("H35").Value = 157
$ws.Range("E37").Value = 201
$ws.Range("F37").Value = 112
$ws.Range("H37").Value = 148
$ws.Range("E38").Value = 108
$ws.Range("F38").Value = 68
$ws.Range("H38").Value = 85
$ws.Range("E40").Value = 318
$ws.Range("F40").Value = 164
$ws.Range("H40").Value = 244
$ws.Range("E41").Value = 451
$ws.Range("F41").Value = 227
$ws.Range("H41").Value = 319
$ws.Range("E42").Value = 482
$ws.Range("F42").Value = 269
$ws.Range("H42").Value = 330
$ws.Range("E43").Value = 149
$ws.Range("F43").Value = 86
$ws.Range("H43").Value = 113
$ws.Range("E44").Value = 387
$ws.Range("E45").Value = 189
$ws.Range("F45").Value = 105
$ws.Range("H45").Value = 144
$ws.Range("E46").Value = 408
$ws.Range("F46").Value = 237
$ws.Range("H46").Value = 301
$ws.Range("E47").Value = 560
$ws.Range("F47").Value = 312
$ws.Range("H47").Value = 404
$ws.Range("E48").Value = 284
$ws.Range("F48").Value = 136
$ws.Range("H48").Value = 180
$ws.Range("E49").Value = 346
$ws.Range("F49").Value = 170
$ws.Range("H49").Value = 257
$ws.Range("E52").Value = 34
